$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Mise en route" / "Tom" hours from 3.5 to 4.5
$ws.Range("B2").Value = 4.5

# Update the active cell selection
$ws.Range("F8").Select()
